# Báo Cáo Thầy Quân
# Underline the "Trần Gia Lương " entries (name + trailing space, plus the
# paragraph mark) in both the member-list table and the task-assignment
# table.

$d = $word.ActiveDocument

for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $table = $d.Tables.Item($i)
    for ($r = 1; $r -le $table.Rows.Count; $r++) {
        $cell = $table.Cell($r, 1)
        $text = $cell.Range.Text
        if ($text -like "*Trần Gia Lương*") {
            $cell.Range.Font.Underline = 1
        }
    }
}
